# Auto-generated edit script: appends the 2025-09-25 01:15 JST Lancers scrape
# to the 'ランサーズ' sheet, re-numbering existing rows and inserting one new
# listing near the top (matches the source commit's diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$rows = @(
    @{ A = '2025-09-25 01:15:28'; B = '【AI活用】データ分析Webサービス開発パートナー募集'; C = 'システム開発'; D = '200,000 円 ~ 300,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399092'; G = 368; H = '🔥AI,Ai ◆開発' }
    @{ A = '2025-09-25 01:15:28'; B = '【SES案件多数】バックエンドエンジニア募集(Java/PHP/Python/Node.js)'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399874'; G = 320; H = '🔥Python ★Java ◆Node.js ○PHP' }
    @{ A = '2025-09-25 01:15:28'; B = 'あなたAIクローン構築パートナー募集・モデル制作&新規依頼'; C = 'システム開発'; D = '100,000 円 ~ 200,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399534'; G = 303; H = '🔥AI,Ai' }
    @{ A = '2025-09-25 01:15:28'; B = 'Excel・Accessベースの改修や追加、Pythonスクレイピングやデータ整形等の開発員募集'; C = 'システム開発'; D = '200,000 円 ~ 300,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399398'; G = 298; H = '🔥Python ◆開発,スクレイピング' }
    @{ A = '2025-09-25 01:15:28'; B = '【急募】カスタマー向けFAQチャットbotの開発依頼'; C = 'システム開発'; D = '1,000,000 円 ~ 3,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399558'; G = 180; H = '★bot ◆開発' }
    @{ A = '2025-09-25 01:15:28'; B = '既存Excelをベースにした短期計画書管理のWebシステム開発'; C = 'システム開発'; D = '100,000 円 ~ 200,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399602'; G = 153; H = '◆開発,システム開発 ◇管理' }
    @{ A = '2025-09-25 01:15:28'; B = '【Flutter+Firebase】社内ポータルアプリ開発のパートナー募集'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399824'; G = 100; H = '◆開発 ◇アプリ' }
    @{ A = '2025-09-25 01:15:28'; B = '【RPA/ブラウザ操作自動化】Webフォーム大量登録の自動化(継続依頼あり)'; C = 'システム開発'; D = '200,000 円 ~ 300,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399631'; G = 88; H = '◆自動化' }
    @{ A = '2025-09-25 01:15:28'; B = 'Googleフォーム × スプレッドシート × GAS 自動化(ストレスチェック診断/台帳保存あり)'; C = 'システム開発'; D = '50,000 円 ~ 100,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399200'; G = 88; H = '◆自動化' }
    @{ A = '2025-09-25 01:15:28'; B = '急募 【急募】Excelで株の保有リストを自動化したいので制作してくださる方募集!'; C = 'システム開発'; D = '20,000 円 ~ 50,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399727'; G = 83; H = '◆自動化' }
    @{ A = '2025-09-25 01:15:28'; B = '【急募】住宅展示場マッチング診断サービスのMVP開発依頼'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399759'; G = 75; H = '◆開発' }
    @{ A = '2025-09-25 01:15:28'; B = '完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします'; C = 'システム開発'; D = '~ 5,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399071'; G = 70; H = '◆効率化' }
    @{ A = '2025-09-25 01:15:28'; B = '【ペットのアバター化】Pawsitiveプロトタイプ開発の依頼'; C = 'システム開発'; D = '200,000 円 ~ 300,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399313'; G = 68; H = '◆開発' }
    @{ A = '2025-09-25 01:15:28'; B = '【急募】PHP・Lalavelでの既存プログラム改修依頼'; C = 'システム開発'; D = '100,000 円 ~ 200,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5396563'; G = 33; H = '○PHP' }
    @{ A = '2025-09-25 01:15:28'; B = '【SES案件多数/リモート可】フルスタックエンジニア募集(フロント〜バック〜クラウドまで歓迎)'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399877'; G = 25; H = $null }
    @{ A = '2025-09-25 01:15:28'; B = '【SES案件多数/リモート可】インフラエンジニア募集(AWS/Linux/NW設計・構築 等歓迎)'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399876'; G = 25; H = $null }
    @{ A = '2025-09-25 01:15:28'; B = '【SES案件多数/リモート可】フロントエンドエンジニア募集(HTML/CSS〜モダンFWまで歓迎)'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399721'; G = 25; H = $null }
    @{ A = '2025-09-25 01:15:28'; B = '〖リモート可〗Delphiエンジニア募集'; C = 'システム開発'; D = '300,000 円 ~ 500,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5341051'; G = 25; H = $null }
    @{ A = '2025-09-25 01:15:28'; B = '【急募】フロントエンドエンジニア募集!(ややWEBコーダー寄り)'; C = 'システム開発'; D = '300,000 円 ~ 500,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399545'; G = 25; H = $null }
    @{ A = '2025-09-25 01:15:28'; B = '【相談から実装まで伴走できる方歓迎】介護・福祉×テクノロジー事例収集の仕組みづくり'; C = 'システム開発'; D = '50,000 円 ~ 100,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5398932'; G = 18; H = $null }
    @{ A = '2025-09-25 01:15:28'; B = '限定公開 PR 限定公開の仕事'; C = 'システム開発'; D = '20,000 円 ~ 50,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399347'; G = 13; H = $null }
    @{ A = '2025-09-25 01:15:28'; B = 'Android kotlin 非同期処理の呼び方'; C = 'システム開発'; D = '5,000 円 ~ 10,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5399765'; G = 10; H = $null }
)

# Clear any rows below the data block that might linger from a larger old sheet,
# then rewrite every data row (2..23) from scratch.
$lastRow = $rows.Count + 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($row.H -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row.H
    } else {
        $ws.Cells.Item($r, 8).ClearContents()
    }
}

# Rebuild the F-column hyperlinks from scratch in row order so the relationship
# ids line up 1:1 with the (re-numbered) rows.
$ws.Hyperlinks.Delete()
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row.F) | Out-Null
}

# Column H widened to fit the longer skill-summary strings.
$ws.Columns.Item(8).ColumnWidth = 28.17

$ws.Range("A1").Select()
